# =====================================================================
# Edit script: applies the commit's changes to BALANÇO_COMPLETO.xlsx
#   1. Delete sheet "UC GERADORA 2"
#   2. "RESUMO " sheet: update F7/G7, F8/G8, clear F9/G9
#   3. "UC GERADORA" sheet: rows 5-16 date / reading / invoice updates
#   4. "UC BENEF. 1" sheet: row 16 date / reading / invoice updates
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Delete the "UC GERADORA 2" sheet entirely
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("UC GERADORA 2").Delete()
$excel.DisplayAlerts = $true

# ---------------------------------------------------------------------
# Helper: write a numeric-looking string into a cell as TEXT while
# preserving the cell's existing style (NumberFormat round-trip avoids
# Excel auto-creating a brand-new style record).
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $origFmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $origFmt
}

# ---------------------------------------------------------------------
# 2. "RESUMO " sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RESUMO ")

Set-TextValue $ws1.Range("F7") "200027499"
$ws1.Range("G7").Value = "RUA CARAI, Q. 12, L. 15, S/N PARQUE ITATIAIA"

Set-TextValue $ws1.Range("F8") "10031998761"
$ws1.Range("G8").Value = "RUA K-4, Q. 01, L. 10, S/N JARDIM ESPLANADA"

$ws1.Range("F9").Value = ""
$ws1.Range("G9").Value = ""

# ---------------------------------------------------------------------
# 3. "UC GERADORA" sheet, rows 5-16
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("UC GERADORA")

$ws2.Range("B5").Value = "19/12/2024"
$ws2.Range("C5").Value = "18/01/2025"
$ws2.Range("I5").Value = 226
$ws2.Range("J5").Value = 181
$ws2.Range("K5").Value = 211
$ws2.Range("N5").Value = 27.39
$ws2.Range("P5").Value = 1418
$ws2.Range("R5").Value = "12378019-5"
$ws2.Range("S5").Value = 16289
$ws2.Range("T5").Value = 16500

$ws2.Range("B6").Value = "18/01/2025"
$ws2.Range("C6").Value = "17/02/2025"
$ws2.Range("I6").Value = 273
$ws2.Range("J6").Value = 160
$ws2.Range("K6").Value = 190
$ws2.Range("N6").Value = 34.19
$ws2.Range("P6").Value = 1531
$ws2.Range("R6").Value = "12378019-5"
$ws2.Range("S6").Value = 16500
$ws2.Range("T6").Value = 16690

# Row 7: B7 used to hold a formula (=C6) and C7 was empty - both become
# literal text dates now.
$ws2.Range("B7").Value = "17/02/2025"
$ws2.Range("C7").Value = "19/03/2025"
$ws2.Range("I7").Value = 204
$ws2.Range("J7").Value = 282
$ws2.Range("K7").Value = 312
$ws2.Range("N7").Value = 52.61
$ws2.Range("P7").Value = 1453
$ws2.Range("R7").Value = "12378019-5"
$ws2.Range("S7").Value = 16690
$ws2.Range("T7").Value = 17002

$ws2.Range("B8").Value = "19/03/2025"
$ws2.Range("C8").Value = "16/04/2025"
$ws2.Range("I8").Value = 155
$ws2.Range("J8").Value = 179
$ws2.Range("K8").Value = 209
$ws2.Range("N8").Value = 44.52
$ws2.Range("P8").Value = 1429
$ws2.Range("R8").Value = "12378019-5"
$ws2.Range("S8").Value = 17002
$ws2.Range("T8").Value = 17211

$ws2.Range("B9").Value = "16/04/2025"
$ws2.Range("C9").Value = "15/05/2025"
$ws2.Range("I9").Value = 159
$ws2.Range("J9").Value = 164
$ws2.Range("K9").Value = 194
$ws2.Range("N9").Value = 44.77
$ws2.Range("P9").Value = 1424
$ws2.Range("R9").Value = "12378019-5"
$ws2.Range("S9").Value = 17211
$ws2.Range("T9").Value = 17405

$ws2.Range("B10").Value = "15/05/2025"
$ws2.Range("C10").Value = "16/06/2025"
$ws2.Range("I10").Value = 154
$ws2.Range("J10").Value = 181
$ws2.Range("K10").Value = 211
$ws2.Range("N10").Value = 48.03
$ws2.Range("P10").Value = 1397
$ws2.Range("R10").Value = "12378019-5"
$ws2.Range("S10").Value = 17405
$ws2.Range("T10").Value = 17616

$ws2.Range("B11").Value = "16/06/2025"
$ws2.Range("C11").Value = "18/07/2025"
$ws2.Range("I11").Value = 252
$ws2.Range("J11").Value = 119
$ws2.Range("K11").Value = 149
$ws2.Range("N11").Value = 42.91
$ws2.Range("P11").Value = 1530
$ws2.Range("R11").Value = "12378019-5"
$ws2.Range("S11").Value = 17616
$ws2.Range("T11").Value = 17765

$ws2.Range("B12").Value = "18/07/2025"
$ws2.Range("C12").Value = "18/08/2025"
$ws2.Range("I12").Value = 147
$ws2.Range("J12").Value = 368
$ws2.Range("K12").Value = 398
$ws2.Range("N12").Value = 46.33
$ws2.Range("P12").Value = 1309
$ws2.Range("R12").Value = "12378019-5"
$ws2.Range("S12").Value = 17765
$ws2.Range("T12").Value = 18163

$ws2.Range("B13").Value = "18/08/2025"
$ws2.Range("C13").Value = "16/09/2025"
$ws2.Range("I13").Value = 83
$ws2.Range("J13").Value = 511
$ws2.Range("K13").Value = 541
$ws2.Range("N13").Value = 76.59
$ws2.Range("P13").Value = 881
$ws2.Range("R13").Value = "12378019-5"
$ws2.Range("S13").Value = 18163
$ws2.Range("T13").Value = 18704

$ws2.Range("B14").Value = "16/09/2025"
$ws2.Range("C14").Value = "16/10/2025"
$ws2.Range("I14").Value = 88
$ws2.Range("J14").Value = 618
$ws2.Range("K14").Value = 648
$ws2.Range("N14").Value = 85.04000000000001
$ws2.Range("P14").Value = 351
$ws2.Range("R14").Value = "12378019-5"
$ws2.Range("S14").Value = 18704
$ws2.Range("T14").Value = 19352

$ws2.Range("B15").Value = "16/10/2025"
$ws2.Range("C15").Value = "18/11/2025"
$ws2.Range("I15").Value = 112
$ws2.Range("J15").Value = 463
$ws2.Range("K15").Value = 618
$ws2.Range("N15").Value = 248.92
$ws2.Range("P15").Value = 0
$ws2.Range("R15").Value = "12378019-5"
$ws2.Range("S15").Value = 19352
$ws2.Range("T15").Value = 19970

$ws2.Range("B16").Value = "18/11/2025"
$ws2.Range("C16").Value = "18/12/2025"
$ws2.Range("I16").Value = 132
$ws2.Range("J16").Value = 132
$ws2.Range("K16").Value = 513
$ws2.Range("N16").Value = 511.7
$ws2.Range("P16").Value = 0
$ws2.Range("R16").Value = "12378019-5"
$ws2.Range("S16").Value = 19970
$ws2.Range("T16").Value = 20483

# ---------------------------------------------------------------------
# 4. "UC BENEF. 1" sheet, row 16
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("UC BENEF. 1")

$ws3.Range("B16").Value = "19/11/2025"
$ws3.Range("C16").Value = "19/12/2025"
$ws3.Range("F16").Value = 481
$ws3.Range("H16").Value = 0
$ws3.Range("J16").Value = 627.24

# K16 becomes a formula; re-apply formats afterwards from a sibling row
# (formula-entry auto-picks a 2-decimal number format otherwise).
$ws3.Range("K16").Formula = "=H16-F16"
$ws3.Range("K15").Copy()
$ws3.Range("K16").PasteSpecial(-4122)

# N16 becomes a formula too; same format fix-up.
$ws3.Range("N16").Formula = "=M16*'UC GERADORA'!O16"
$ws3.Range("N15").Copy()
$ws3.Range("N16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# R16 is removed entirely (cleared).
$ws3.Range("R16").Value = ""

$ws3.Range("S16").Value = "12340952-7"
$ws3.Range("T16").Value = 11694
$ws3.Range("U16").Value = 12175
